# Fix three small typos in the contract text:
#  1) "Post_condizioni" -> "Post-condizioni" (Contratto CO2 postconditions)
#  2) "Post_condizioni" -> "Post-condizioni" (Contratto CO4 postconditions)
#  3) "nuova istanza nuova di Risorsa" -> "nuova istanza nuovar di Risorsa" (Contratto CO3 postconditions)
#
# (The rest of the upstream diff is just Word's own run-splitting around
# proofing-error markers, i.e. spelling/grammar squiggles, collapsing back
# into single runs with the same visible text -- no content change there.)

$d = $word.ActiveDocument

# wdReplaceAll = 2
$d.Content.Find.Execute(
    "Post_condizioni",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Post-condizioni",
    2
)

$d.Content.Find.Execute(
    "nuova istanza nuova di Risorsa",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "nuova istanza nuovar di Risorsa",
    2
)
